# Applies the cryptos list update described in the commit diff.
# Pattern for numeric-looking price strings (column D): force text storage by
# setting NumberFormat to Text ("@") before assignment, then ClearFormats() to
# drop the temporary number-format style again (Excel would otherwise silently
# parse strings like "292.00" or "1.00" into the numbers 292 / 1, dropping the
# trailing zero that the source data relies on).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "39.945.80"
$ws.Range("D2").ClearFormats()

$ws.Range("E2").Value = "  +0.23%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.216.45"
$ws.Range("D3").ClearFormats()

$ws.Range("E3").Value = "  -0.18%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "292.00"
$ws.Range("D5").ClearFormats()

$ws.Range("E5").Value = "  -0.07%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "86.84"
$ws.Range("D6").ClearFormats()

$ws.Range("E6").Value = "  -0.72%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.515"
$ws.Range("D7").ClearFormats()

$ws.Range("E7").Value = "  -0.35%  "

$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.467"
$ws.Range("D9").ClearFormats()

$ws.Range("E9").Value = "  -1.12%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "30.42"
$ws.Range("D10").ClearFormats()

$ws.Range("E10").Value = "  +0.17%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "50.28"
$ws.Range("D11").ClearFormats()

$ws.Range("E11").Value = "  +5.88%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0778"
$ws.Range("D12").ClearFormats()

$ws.Range("E12").Value = "  -1.01%  "

$ws.Range("E13").Value = "  +3.10%  "

$ws.Range("E14").Value = "  +0.40%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.557.96"
$ws.Range("D15").ClearFormats()

$ws.Range("E15").Value = "  -0.25%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.75"
$ws.Range("D16").ClearFormats()

$ws.Range("E16").Value = "  -2.49%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.181.92"
$ws.Range("D17").ClearFormats()

$ws.Range("E17").Value = "  -1.50%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.731"
$ws.Range("D18").ClearFormats()

$ws.Range("E18").Value = "  +0.08%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "39.886.73"
$ws.Range("D19").ClearFormats()

$ws.Range("E19").Value = "  +0.18%  "

$ws.Range("E20").Value = "  +0.21%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.11"
$ws.Range("D21").ClearFormats()

$ws.Range("E21").Value = "  -3.65%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.74"
$ws.Range("D22").ClearFormats()

$ws.Range("E22").Value = "  -1.53%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "65.59"
$ws.Range("D23").ClearFormats()

$ws.Range("E23").Value = "  -0.27%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "237.49"
$ws.Range("D24").ClearFormats()

$ws.Range("E24").Value = "  +0.73%  "

$ws.Range("E25").Value = "  +0.04%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.46"
$ws.Range("D26").ClearFormats()

$ws.Range("E26").Value = "  -0.55%  "

$ws.Range("E27").Value = "  -0.07%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "23.09"
$ws.Range("D28").ClearFormats()

$ws.Range("E28").Value = "  +1.13%  "

$ws.Range("B29").Value = "Toncoin"

$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.14"
$ws.Range("D29").ClearFormats()

$ws.Range("E29").Value = "  +0.29%  "

$ws.Range("B30").Value = "Cosmos"

$ws.Range("C30").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.23"
$ws.Range("D30").ClearFormats()

$ws.Range("E30").Value = "  -0.33%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "156.85"
$ws.Range("D31").ClearFormats()

$ws.Range("E31").Value = "  +2.73%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "31.83"
$ws.Range("D32").ClearFormats()

$ws.Range("E32").Value = "  -2.85%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.00"
$ws.Range("D33").ClearFormats()

$ws.Range("E33").Value = "  +0.07%  "

$ws.Range("E34").Value = "  +0.20%  "

$ws.Range("E35").Value = "  +5.94%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0713"
$ws.Range("D36").ClearFormats()

$ws.Range("E36").Value = "  -0.84%  "

$ws.Range("E37").Value = "  -1.89%  "

$ws.Range("E38").Value = "  -0.30%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0990"
$ws.Range("D39").ClearFormats()

$ws.Range("E39").Value = "  -0.17%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.72"
$ws.Range("D40").ClearFormats()

$ws.Range("E40").Value = "  +0.64%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "15.19"
$ws.Range("D41").ClearFormats()

$ws.Range("E41").Value = "  -4.66%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.091.05"
$ws.Range("D42").ClearFormats()

$ws.Range("E42").Value = "  -0.24%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.70"
$ws.Range("D43").ClearFormats()

$ws.Range("E43").Value = "  -2.67%  "

$ws.Range("E44").Value = "  +0.61%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "17.93"
$ws.Range("D45").ClearFormats()

$ws.Range("E45").Value = "  +1.05%  "

$ws.Range("E46").Value = "  -2.11%  "

$ws.Range("E47").Value = "  -8.12%  "

$ws.Range("E48").Value = "  +2.43%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.431.62"
$ws.Range("D49").ClearFormats()

$ws.Range("E49").Value = "  -0.18%  "

$ws.Range("E50").Value = "  +0.18%  "

$ws.Range("E51").Value = "  +2.60%  "
